$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header cells, matching the formatting used by the other header cells
# (bold font, thin box border, centered/top aligned - same as e.g. AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous

# Fill in the win/loss/tie record for every data row (2 through 51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 64   # AD
    $ws.Cells.Item($r, 31).Value = 98   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
